$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - (Intercept)
$ws.Range("B2").Value = 9405.951567
$ws.Range("D2").Value = 125.667214

# Row 3 - household_group_collapsed
$ws.Range("B3").Value = 1441.540078
$ws.Range("D3").Value = 9.629771
$ws.Range("E3").Value = 0.000097

# Row 4 - Residuals
$ws.Range("B4").Value = 16691.125237
$ws.Range("C4").Value = 223

# Row 5 - SM-Control
$ws.Range("G5").Value = -5.090259
$ws.Range("H5").Value = -8.868672999999999
$ws.Range("I5").Value = -1.311846
$ws.Range("J5").Value = 0.004793

# Row 6 - SM + Traps-Control
$ws.Range("G6").Value = -0.058084
$ws.Range("H6").Value = -4.126635
$ws.Range("I6").Value = 4.010468
$ws.Range("J6").Value = 0.999375

# Row 7 - SM + Traps-SM
$ws.Range("G7").Value = 5.032176
$ws.Range("H7").Value = 1.956309
$ws.Range("I7").Value = 8.108043
$ws.Range("J7").Value = 0.000435
